$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7 ("Polynomial Regression"), shifting rows 7-11 down to 8-12.
$ws.Rows.Item(7).Insert()

# Copy formatting from the row below (now row 8, originally row 7) into the new row 7
# so the new row inherits the same body-row style (borders etc.).
$ws.Range("A8:C8").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)

# Fill in the new row's data: Id, Model name, Accuracy value.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Lasso Regression+normalization+lag1+PCA(2)"
$ws.Range("C7").Value = 73.409172663976904

# Renumber the Id column for all rows pushed down below the new row.
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10

# Widen column B to fit the longer text.
$ws.Columns.Item(2).ColumnWidth = 40.5

# Match the selection shown in the saved workbook.
$ws.Range("B7").Select()
